# Update the cached "today" date shown in every date placeholder
# (slide master, all slide layouts, notes master, and every slide)
# from 2020/10/7 to 2020/10/8, and retitle the "10 Query: Part 1"
# title text to "10 Query Custom Field: Part 1" on the three slides
# that show it.

$p = $ppt.ActivePresentation

$oldDate = "2020/10/7"
$newDate = "2020/10/8"

# --- Slide master: date placeholder is shape #3 ---
$p.SlideMaster.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# --- Slide layouts: date placeholder shape index per layout ---
$layoutDateIndex = @{
    1  = 3
    2  = 3
    3  = 3
    4  = 4
    5  = 6
    6  = 2
    7  = 1
    8  = 4
    9  = 4
    10 = 3
    11 = 3
}

$layouts = $p.SlideMaster.CustomLayouts
foreach ($key in $layoutDateIndex.Keys) {
    $layout = $layouts.Item($key)
    $layout.Shapes.Item($layoutDateIndex[$key]).TextFrame.TextRange.Text = $newDate
}

# --- Notes master: date placeholder is shape #2 ---
$p.NotesMaster.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# --- Every slide: date placeholder shape index per slide ---
$slideDateIndex = @{
    1  = 3
    2  = 4
    3  = 4
    4  = 2
    5  = 4
    6  = 2
    7  = 4
    8  = 2
    9  = 4
    10 = 2
    11 = 4
    12 = 2
    13 = 4
    14 = 2
    15 = 4
    16 = 4
    17 = 4
    18 = 2
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $slide.Shapes.Item($slideDateIndex[$i]).TextFrame.TextRange.Text = $newDate
}

# --- Title text update on slides 1-3: title is always shape #1 ---
$oldTitle = "10 Query: Part 1"
$newTitle = "10 Query Custom Field: Part 1"

foreach ($idx in 1..3) {
    $slide = $p.Slides.Item($idx)
    $titleShape = $slide.Shapes.Item(1)
    if ($titleShape.TextFrame.TextRange.Text -eq $oldTitle) {
        $titleShape.TextFrame.TextRange.Text = $newTitle
    }
}
